$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '30.113.42'
$cell.Style = $origStyle
$cell = $ws.Range("E2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +5.60%  '
$cell.Style = $origStyle
$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.919.60'
$cell.Style = $origStyle
$cell = $ws.Range("E3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.47%  '
$cell.Style = $origStyle
$cell = $ws.Range("D4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = $origStyle
$cell = $ws.Range("E4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.56%  '
$cell.Style = $origStyle
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '330.28'
$cell.Style = $origStyle
$cell = $ws.Range("E5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.60%  '
$cell.Style = $origStyle
$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = $origStyle
$cell = $ws.Range("E6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.57%  '
$cell.Style = $origStyle
$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.5213'
$cell.Style = $origStyle
$cell = $ws.Range("E7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.32%  '
$cell.Style = $origStyle
$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.4084'
$cell.Style = $origStyle
$cell = $ws.Range("E8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.58%  '
$cell.Style = $origStyle
$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.08520'
$cell.Style = $origStyle
$cell = $ws.Range("E9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.89%  '
$cell.Style = $origStyle
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '43.32'
$cell.Style = $origStyle
$cell = $ws.Range("E10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.70%  '
$cell.Style = $origStyle
$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.128'
$cell.Style = $origStyle
$cell = $ws.Range("E11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.09%  '
$cell.Style = $origStyle
$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '22.37'
$cell.Style = $origStyle
$cell = $ws.Range("E12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +9.58%  '
$cell.Style = $origStyle
$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.409'
$cell.Style = $origStyle
$cell = $ws.Range("E13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.01%  '
$cell.Style = $origStyle
$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.922.69'
$cell.Style = $origStyle
$cell = $ws.Range("E14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.19%  '
$cell.Style = $origStyle
$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.415'
$cell.Style = $origStyle
$cell = $ws.Range("E15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.75%  '
$cell.Style = $origStyle
$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = $origStyle
$cell = $ws.Range("E16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.55%  '
$cell.Style = $origStyle
$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '95.68'
$cell.Style = $origStyle
$cell = $ws.Range("E17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.87%  '
$cell.Style = $origStyle
$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.00001113'
$cell.Style = $origStyle
$cell = $ws.Range("E18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.81%  '
$cell.Style = $origStyle
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.06715'
$cell.Style = $origStyle
$cell = $ws.Range("E19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.26%  '
$cell.Style = $origStyle
$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '18.31'
$cell.Style = $origStyle
$cell = $ws.Range("E20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.21%  '
$cell.Style = $origStyle
$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = $origStyle
$cell = $ws.Range("E21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.59%  '
$cell.Style = $origStyle
$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.032'
$cell.Style = $origStyle
$cell = $ws.Range("E22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.03%  '
$cell.Style = $origStyle
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '30.126.86'
$cell.Style = $origStyle
$cell = $ws.Range("E23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +5.50%  '
$cell.Style = $origStyle
$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '11.36'
$cell.Style = $origStyle
$cell = $ws.Range("E24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.89%  '
$cell.Style = $origStyle
$cell = $ws.Range("E25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.03%  '
$cell.Style = $origStyle
$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.144.70'
$cell.Style = $origStyle
$cell = $ws.Range("E26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.50%  '
$cell.Style = $origStyle
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '21.13'
$cell.Style = $origStyle
$cell = $ws.Range("E27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.33%  '
$cell.Style = $origStyle
$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '160.01'
$cell.Style = $origStyle
$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.451'
$cell.Style = $origStyle
$cell = $ws.Range("E29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.54%  '
$cell.Style = $origStyle
$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '129.10'
$cell.Style = $origStyle
$cell = $ws.Range("E30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.03%  '
$cell.Style = $origStyle
$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.079'
$cell.Style = $origStyle
$cell = $ws.Range("E31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.87%  '
$cell.Style = $origStyle
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.1056'
$cell.Style = $origStyle
$cell = $ws.Range("E32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.41%  '
$cell.Style = $origStyle
$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.087'
$cell.Style = $origStyle
$cell = $ws.Range("E33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +6.10%  '
$cell.Style = $origStyle
$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.640'
$cell.Style = $origStyle
$cell = $ws.Range("E34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.79%  '
$cell.Style = $origStyle
$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.02501'
$cell.Style = $origStyle
$cell = $ws.Range("E35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.05%  '
$cell.Style = $origStyle
$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.06606'
$cell.Style = $origStyle
$cell = $ws.Range("E36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.78%  '
$cell.Style = $origStyle
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.2217'
$cell.Style = $origStyle
$cell = $ws.Range("E37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.45%  '
$cell.Style = $origStyle
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.217'
$cell.Style = $origStyle
$cell = $ws.Range("E38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.90%  '
$cell.Style = $origStyle
$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.234'
$cell.Style = $origStyle
$cell = $ws.Range("E39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.65%  '
$cell.Style = $origStyle
$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.944'
$cell.Style = $origStyle
$cell = $ws.Range("E40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.03%  '
$cell.Style = $origStyle
$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.6545'
$cell.Style = $origStyle
$cell = $ws.Range("E41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.56%  '
$cell.Style = $origStyle
$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '11.70'
$cell.Style = $origStyle
$cell = $ws.Range("E42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +5.45%  '
$cell.Style = $origStyle
$cell = $ws.Range("E43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.13%  '
$cell.Style = $origStyle
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.6168'
$cell.Style = $origStyle
$cell = $ws.Range("E44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.75%  '
$cell.Style = $origStyle
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '13.33'
$cell.Style = $origStyle
$cell = $ws.Range("E45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.57%  '
$cell.Style = $origStyle
$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.774'
$cell.Style = $origStyle
$cell = $ws.Range("E46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.32%  '
$cell.Style = $origStyle
$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.088'
$cell.Style = $origStyle
$cell = $ws.Range("E47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.35%  '
$cell.Style = $origStyle
$cell = $ws.Range("E48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.87%  '
$cell.Style = $origStyle
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '124.94'
$cell.Style = $origStyle
$cell = $ws.Range("E49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.49%  '
$cell.Style = $origStyle
$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.163'
$cell.Style = $origStyle
$cell = $ws.Range("E50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.57%  '
$cell.Style = $origStyle
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '79.77'
$cell.Style = $origStyle
$cell = $ws.Range("E51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.32%  '
$cell.Style = $origStyle
